# ---------------------------------------------------------------------------
# Word-COM-interop script implementing the "cleaned up related works" commit.
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# PART 1: Related-works paragraph reshuffle (WPA-cracking section)
# ---------------------------------------------------------------------------

# 1a. Drop the stray "_GoBack" bookmark that currently sits at the end of the
#     "...Chopchop attack [12]." paragraph -- it gets re-homed later.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 1b. Delete the old paragraphs in between "[12]." and the "SECURING WIFI"
#     heading: the GPGPU/Pyrit paragraph, "Many have abandoned WEP...",
#     "Possibly the first attack...", the TODO paragraph and the trailing
#     empty paragraph.
$rStart = $d.Content.Duplicate
$rStart.Find.Execute("As the age of general purpose graphics processing units", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$rEnd = $d.Content.Duplicate
$rEnd.Find.Execute("SECURING WIFI", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $rEnd.Paragraphs.First

$oldBlock = $d.Range($rStart.Start, $headingPara.Range.Start)
$oldBlock.Delete()

# 1c. Insert the new paragraphs (in their new order/wording) right after the
#     "...Chopchop attack [12]." paragraph (i.e. right before "SECURING
#     WIFI"). Inserting at (paragraph-end - 1) -- i.e. just before that
#     paragraph's own mark -- reliably splices new paragraphs in after it
#     without clobbering neighbouring content.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Chopchop attack", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchor.Paragraphs.First
$insertionPoint = $d.Range($anchorPara.Range.End - 1, $anchorPara.Range.End - 1)

$newXml = "<w:p $wNs>" +
    "<w:r><w:lastRenderedPageBreak/><w:t>Many have abandoned WEP, saying that it is &#8220;completely insecure&#8221;</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> [8]</w:t></w:r>" +
    "<w:r><w:t>.  WPA was its replacement.  Eventually, WPA was replaced with WPA2.  So far these have remained relatively secure.  The same number of vulnerabilities in WEP has not been found in WPA and WPA2.  A number of attacks do exist though.</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
    "<w:r><w:t>Possibly the first attack on WPA was presented in [</w:t></w:r>" +
    "<w:r><w:t>6</w:t></w:r>" +
    "<w:r><w:t>].</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>  It demonstrated that a </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>chopchop</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> like attack (an attack used to crack WEP) could also be used to crack WPA.  A variety of conditions must be met for this attack to work on WPA.  However, it was shown that these conditions are not unreasonable in most wireless networks [6].  </w:t></w:r>" +
    "<w:r><w:t>Additionally, this attack only works for WPA with TKIP not CCMP.</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
    "<w:r><w:t xml:space='preserve'>As the age of general purpose graphics processing units (GPGPU) computing dawned, </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>the extreme parallelism offered by the GPU became </w:t></w:r>" +
    "<w:r><w:t>clearer</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'> and more widely used.</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>  GPGPU computing is now being used to crack WPA encrypted wireless networks.  </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>Pyrit</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> is one such approach [9].</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>  Using CUDA [</w:t></w:r>" +
    "<w:r><w:t>13</w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>], </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>Pyrit</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> can compute up to 89,000 pairwise master keys per second.</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wNs>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
    "<w:bookmarkEnd w:id='0'/>" +
    "</w:p>"

$insertionPoint.InsertXML($newXml)

Write-Output "Part 1 complete"
